$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "3：12-6：21"
$ws.Range("C18").Value = "函数递归 写了些代码，还改进了些代码"
$ws.Range("E18").Value = "（文档里的不是很全，后面有些代码我就直接记录的，没有自己打）"

$ws.Range("E18").Select()
